$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24
$url = "https://www.genomeweb.com/cancer/foresight-diagnostics-launches-trial-mrd-guided-treatment-hodgkin-lymphoma"
$keyword = "CDx, ctDNA"
$title = "Foresight Diagnostics Launches Trial of MRD-Guided Treatment in Hodgkin Lymphoma"

$ws.Range("B$newRow").Value = $keyword
$ws.Range("C$newRow").Value = $title

$ws.Hyperlinks.Add($ws.Range("A$newRow"), $url, "", "", $url)
